$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 4 data rows (rows 2-5, corresponding to years 1985-1988).
# This shifts all remaining rows up by 4, so old row 6 becomes new row 2, etc.
$ws.Range("A2:E5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
